$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert two new quarter columns (D:E) in front of the existing data, which
# shifts the previous D:K block to F:M. This mirrors pasting two new
# columns of quarterly financials into the three statements (Income
# Statement, Balance Sheet, Cash Flow Statement) that all share the same
# column layout, each starting at row 7, 38 and 80 respectively.
# ---------------------------------------------------------------------------
$ws.Columns("D:E").Insert()

# Copy the number/date formatting from the (now shifted) first data column
# pair F:G onto the newly inserted, still-blank D:E columns so every row
# keeps its original style (date style for the "Period Ending" rows, the
# right-aligned #,##0 style for every other data row).
$ws.Range("F7:G102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Fill in the actual values for the two new quarters (D = newest quarter,
# E = the quarter right before it) for every row that carries data.
# ---------------------------------------------------------------------------
$rows  = @(7,8,9,10,12,13,14,15,17,18,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,38,41,42,43,44,45,46,47,48,49,50,51,52,53,54,57,58,59,60,61,62,63,64,65,66,68,69,70,71,72,73,74,75,76,77,80,81,83,84,85,86,87,88,89,91,92,93,94,96,97,98,99,100,101,102)
$dvals = @(43465,2100,1100,1000,200,0,0,0,2200,-100,0,0,0,-100,0,0,-100,-100,0,0,0,0,0,-100,0,-100,43465,600,0,1000,1400,200,3100,0,200,300,0,0,0,0,3700,700,"NA",400,900,0,100,0,0,0,1100,0,0,0,0,-21700,0,0,0,2500,0,43465,-100,0,0,0,0,0,0,100,0,0,0,0,0,0,0,0,400,0,400)
$evals = @(43373,2200,1000,1200,200,0,0,0,2200,0,0,100,0,0,0,0,0,0,0,0,0,0,0,0,0,0,43373,200,0,1100,1400,100,2800,0,300,300,0,0,0,0,3300,600,"NA",400,900,0,0,0,0,0,1100,0,0,0,0,-21600,0,0,0,2300,0,43373,0,0,0,0,0,0,0,200,0,0,0,0,0,0,0,0,-100,0,100)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $ws.Cells.Item($r, 4).Value = $dvals[$i]
    $ws.Cells.Item($r, 5).Value = $evals[$i]
}
